$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.065182100734830328
$ws.Range("B1").Value = 0.065112778237065072
$ws.Range("A2").Value = -0.019002471753216454
$ws.Range("B2").Value = 0.018830402752247011
$ws.Range("A3").Value = 0.084100224335667662
$ws.Range("B3").Value = -0.0842725142130476
$ws.Range("A4").Value = -0.19972060499146593
$ws.Range("B4").Value = 0.19876749690787321
$ws.Range("A5").Value = -0.19276749740385224
$ws.Range("B5").Value = 0.19084084150214053
$ws.Range("A6").Value = -0.075953957097774083
$ws.Range("B6").Value = 0.075884224100142728
$ws.Range("A7").Value = -0.055884224689277673
$ws.Range("B7").Value = 0.055732213795613106
$ws.Range("A8").Value = -0.035732214388622729
$ws.Range("B8").Value = 0.035607157562709268
$ws.Range("A9").Value = -0.029607158084484553
$ws.Range("B9").Value = 0.029499996731691347
$ws.Range("A10").Value = -0.056010284887200612
$ws.Range("B10").Value = 0.055964901731492489
$ws.Range("A11").Value = -0.051464902249300337
$ws.Range("B11").Value = 0.051387858519429841
$ws.Range("A12").Value = -0.045387859046467582
$ws.Range("B12").Value = 0.045148411408242328
$ws.Range("A13").Value = -0.039148411942609762
$ws.Range("B13").Value = 0.039083245602221872
$ws.Range("A14").Value = -0.02708324617087321
$ws.Range("B14").Value = 0.02705179519013079
$ws.Range("A15").Value = -0.021051795728099343
$ws.Range("B15").Value = 0.021027045427256752
$ws.Range("A16").Value = -0.015027045966633512
$ws.Range("B16").Value = 0.015004350236447639
$ws.Range("A17").Value = -0.0090043507777037846
$ws.Range("B17").Value = 0.0089999994423592966
$ws.Range("A18").Value = -0.036110306994746821
$ws.Range("B18").Value = 0.036096539942715111
$ws.Range("A19").Value = -0.027096540448479534
$ws.Range("B19").Value = 0.027013683443542824
$ws.Range("A20").Value = -0.018013683953347126
$ws.Range("B20").Value = 0.01800427928588455
$ws.Range("A21").Value = -0.0090042797962324173
$ws.Range("B21").Value = 0.0089999994893625868
$ws.Range("A22").Value = -0.093930627594701832
$ws.Range("B22").Value = 0.093623188490546738
$ws.Range("A23").Value = -0.084623189001008292
$ws.Range("B23").Value = 0.084124552113959083
$ws.Range("A24").Value = -0.042124552809021942
$ws.Range("B24").Value = 0.041999999301569169
$ws.Range("A25").Value = -0.1088868849170197
$ws.Range("B25").Value = 0.10873369492106377
$ws.Range("A26").Value = -0.10273369543532951
$ws.Range("B26").Value = 0.10254090181367559
$ws.Range("A27").Value = -0.096540902330200851
$ws.Range("B27").Value = 0.095897620873736322
$ws.Range("A28").Value = -0.089897621401532568
$ws.Range("B28").Value = 0.089469413264342101
$ws.Range("A29").Value = -0.066678005559190012
$ws.Range("B29").Value = 0.066420989862747248
$ws.Range("A30").Value = -0.042164117318785888
$ws.Range("B30").Value = 0.04201920206338805
$ws.Range("A31").Value = -0.027019202655111485
$ws.Range("B31").Value = 0.027000796084028167
$ws.Range("A32").Value = -0.0060007967085908831
$ws.Range("B32").Value = 0.0059999994557982106
